# Updates cryptos list data cells (Price column D, Volume(1h) column E, and
# for a few rows that were reordered/swapped also Coin (B) and Link (C)).
# Values are written as literal text (not auto-converted to numbers) by
# temporarily forcing a text number format, then restoring the "Normal"
# style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"; $ws.Range("D2").Value = '90.822.53'; $ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"; $ws.Range("E2").Value = '  -0.46%  '; $ws.Range("E2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"; $ws.Range("D3").Value = '3.152.47'; $ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"; $ws.Range("E3").Value = '  +2.19%  '; $ws.Range("E3").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"; $ws.Range("E4").Value = '  +0.27%  '; $ws.Range("E4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"; $ws.Range("D5").Value = '215.66'; $ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"; $ws.Range("E5").Value = '  -0.05%  '; $ws.Range("E5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"; $ws.Range("D6").Value = '625.49'; $ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"; $ws.Range("E6").Value = '  +1.21%  '; $ws.Range("E6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"; $ws.Range("D7").Value = '1.17'; $ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"; $ws.Range("E7").Value = '  +33.72%  '; $ws.Range("E7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"; $ws.Range("D8").Value = '0.367'; $ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"; $ws.Range("E8").Value = '  -2.03%  '; $ws.Range("E8").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"; $ws.Range("E9").Value = '  -0.03%  '; $ws.Range("E9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"; $ws.Range("D10").Value = '3.149.59'; $ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"; $ws.Range("E10").Value = '  +2.25%  '; $ws.Range("E10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"; $ws.Range("D11").Value = '0.757'; $ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"; $ws.Range("E11").Value = '  +13.17%  '; $ws.Range("E11").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"; $ws.Range("E12").Value = '  +6.77%  '; $ws.Range("E12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"; $ws.Range("D13").Value = '5.74'; $ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"; $ws.Range("E13").Value = '  +7.09%  '; $ws.Range("E13").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"; $ws.Range("E14").Value = '  -0.65%  '; $ws.Range("E14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"; $ws.Range("D15").Value = '35.01'; $ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"; $ws.Range("E15").Value = '  +6.57%  '; $ws.Range("E15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"; $ws.Range("D16").Value = '90.582.08'; $ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"; $ws.Range("D17").Value = '3.737.66'; $ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"; $ws.Range("E17").Value = '  +2.45%  '; $ws.Range("E17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"; $ws.Range("D18").Value = '3.192.93'; $ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"; $ws.Range("E18").Value = '  +3.71%  '; $ws.Range("E18").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"; $ws.Range("E19").Value = '  +7.35%  '; $ws.Range("E19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"; $ws.Range("D20").Value = '14.64'; $ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"; $ws.Range("E20").Value = '  +6.55%  '; $ws.Range("E20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"; $ws.Range("D21").Value = '473.06'; $ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"; $ws.Range("E21").Value = '  +9.15%  '; $ws.Range("E21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"; $ws.Range("D22").Value = '0.0000210'; $ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"; $ws.Range("E22").Value = '  -4.94%  '; $ws.Range("E22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"; $ws.Range("D23").Value = '9.14'; $ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"; $ws.Range("E23").Value = '  +8.12%  '; $ws.Range("E23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"; $ws.Range("D24").Value = '5.33'; $ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"; $ws.Range("E24").Value = '  +4.75%  '; $ws.Range("E24").Style = "Normal"
$ws.Range("B25").NumberFormat = "@"; $ws.Range("B25").Value = 'Litecoin'; $ws.Range("B25").Style = "Normal"
$ws.Range("C25").NumberFormat = "@"; $ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; $ws.Range("C25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"; $ws.Range("D25").Value = '94.97'; $ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"; $ws.Range("E25").Value = '  +13.47%  '; $ws.Range("E25").Style = "Normal"
$ws.Range("B26").NumberFormat = "@"; $ws.Range("B26").Value = 'NEARProtocol'; $ws.Range("B26").Style = "Normal"
$ws.Range("C26").NumberFormat = "@"; $ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; $ws.Range("C26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"; $ws.Range("D26").Value = '5.75'; $ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"; $ws.Range("E26").Value = '  +3.99%  '; $ws.Range("E26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"; $ws.Range("D27").Value = '12.38'; $ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"; $ws.Range("E27").Value = '  +4.44%  '; $ws.Range("E27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"; $ws.Range("D28").Value = '3.327.77'; $ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"; $ws.Range("E28").Value = '  +2.95%  '; $ws.Range("E28").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"; $ws.Range("E29").Value = '  -0.15%  '; $ws.Range("E29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"; $ws.Range("D30").Value = '9.32'; $ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"; $ws.Range("E30").Value = '  +8.14%  '; $ws.Range("E30").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"; $ws.Range("E31").Value = '  -3.08%  '; $ws.Range("E31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"; $ws.Range("D32").Value = '0.214'; $ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"; $ws.Range("E32").Value = '  +53.69%  '; $ws.Range("E32").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"; $ws.Range("E33").Value = '  -7.23%  '; $ws.Range("E33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"; $ws.Range("D34").Value = '27.46'; $ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"; $ws.Range("E34").Value = '  +19.48%  '; $ws.Range("E34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"; $ws.Range("D35").Value = '518.31'; $ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"; $ws.Range("E35").Value = '  +0.71%  '; $ws.Range("E35").Style = "Normal"
$ws.Range("B36").NumberFormat = "@"; $ws.Range("B36").Value = 'Kaspa'; $ws.Range("B36").Style = "Normal"
$ws.Range("C36").NumberFormat = "@"; $ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; $ws.Range("C36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"; $ws.Range("D36").Value = '0.146'; $ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"; $ws.Range("E36").Value = '  +6.92%  '; $ws.Range("E36").Style = "Normal"
$ws.Range("B37").NumberFormat = "@"; $ws.Range("B37").Value = 'PancakeSwap'; $ws.Range("B37").Style = "Normal"
$ws.Range("C37").NumberFormat = "@"; $ws.Range("C37").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; $ws.Range("C37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"; $ws.Range("D37").Value = '1.94'; $ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"; $ws.Range("E37").Value = '  +5.59%  '; $ws.Range("E37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"; $ws.Range("D38").Value = '3.59'; $ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"; $ws.Range("E38").Value = '  -6.22%  '; $ws.Range("E38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"; $ws.Range("D39").Value = '6.93'; $ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"; $ws.Range("E39").Value = '  +1.09%  '; $ws.Range("E39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"; $ws.Range("D40").Value = '1.30'; $ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"; $ws.Range("E40").Value = '  +3.03%  '; $ws.Range("E40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"; $ws.Range("D41").Value = '0.0915'; $ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"; $ws.Range("E41").Value = '  +27.53%  '; $ws.Range("E41").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"; $ws.Range("E42").Value = '  -0.41%  '; $ws.Range("E42").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"; $ws.Range("E43").Value = '  +15.81%  '; $ws.Range("E43").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"; $ws.Range("E44").Value = '  +0.08%  '; $ws.Range("E44").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"; $ws.Range("E45").Value = '  +6.29%  '; $ws.Range("E45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"; $ws.Range("D46").Value = '0.747'; $ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"; $ws.Range("E46").Value = '  +23.73%  '; $ws.Range("E46").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"; $ws.Range("E47").Value = '  +0.00%  '; $ws.Range("E47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"; $ws.Range("D48").Value = '4.69'; $ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"; $ws.Range("E48").Value = '  +10.81%  '; $ws.Range("E48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"; $ws.Range("D49").Value = '150.80'; $ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"; $ws.Range("E49").Value = '  +5.34%  '; $ws.Range("E49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"; $ws.Range("D50").Value = '45.51'; $ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"; $ws.Range("E50").Value = '  +4.23%  '; $ws.Range("E50").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"; $ws.Range("E51").Value = '  +9.64%  '; $ws.Range("E51").Style = "Normal"
